# 7.10 Fixed Some Bugs
# Wrap the three "note comparison" dialogue lines in a green color tag
# (Unity-style rich text) and adjust the row heights of rows 2 and 3 to
# accommodate the now-longer wrapped text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = " <color=#00CC00>(If the note had been written under normal circumstances, what would it typically look like?)</color>"
$ws.Range("B3").Value = " <color=#00CC00>(Now compare that to the current note.)</color>"
$ws.Range("B4").Value = " <color=#00CC00>(Don’t you think there’s something a bit strange about it?)</color>"

$ws.Rows.Item(2).RowHeight = 51
$ws.Rows.Item(3).RowHeight = 34
